# Dagbaekur fyrir alvoru karlmenn
# Add hours logged on Wednesday (column E) for week-4 block (rows 26-32):
#   row27 (Kröfulýsing)  -> E27 = 180
#   row28 (Hönnun)       -> E28 = 60
#   row29 (Rannsóknir)   -> E29 = 60
# Everything downstream (J27:J29, J32, D57:D59, D62) is formula-driven and
# recalculates automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

$ws.Range("E27").Value = 180
$ws.Range("E28").Value = 60
$ws.Range("E29").Value = 60

# Restore the view/selection state recorded for this sheet.
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("E30").Select()
